$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.126.66'
$ws.Range("E2").Value = '  +8.59%  '
$ws.Range("D3").Value = '3.450.97'
$ws.Range("E3").Value = '  +5.65%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'414.62"
$ws.Range("E5").Value = '  +4.07%  '
$ws.Range("D6").Value = "'125.59"
$ws.Range("E6").Value = '  +15.36%  '
$ws.Range("D7").Value = '3.445.36'
$ws.Range("E7").Value = '  +5.69%  '
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +6.14%  '
$ws.Range("E11").Value = '  +32.95%  '
$ws.Range("D12").Value = "'41.60"
$ws.Range("E12").Value = '  +5.56%  '
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = '3.991.31'
$ws.Range("E14").Value = '  +5.68%  '
$ws.Range("D15").Value = "'8.52"
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").Value = "'19.78"
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").Value = '3.443.42'
$ws.Range("E17").Value = '  +5.21%  '
$ws.Range("D18").Value = '62.070.06'
$ws.Range("E18").Value = '  +8.86%  '
$ws.Range("D19").Value = "'1.04"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = "'11.03"
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = "'0.0000132"
$ws.Range("E21").Value = '  +20.73%  '
$ws.Range("D22").Value = "'3.34"
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = "'82.28"
$ws.Range("E23").Value = '  +10.81%  '
$ws.Range("D24").Value = "'315.17"
$ws.Range("E24").Value = '  +7.09%  '
$ws.Range("D25").Value = "'13.00"
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  +10.38%  '
$ws.Range("D28").Value = "'7.86"
$ws.Range("E28").Value = '  +5.89%  '
$ws.Range("D29").Value = "'7.90"
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  +2.55%  '
$ws.Range("E32").Value = '  +4.27%  '
$ws.Range("D33").Value = "'11.58"
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = '  +19.93%  '
$ws.Range("D35").Value = "'42.14"
$ws.Range("E35").Value = '  +4.87%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").Value = "'52.25"
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("D39").Value = "'3.51"
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("D42").Value = "'2.01"
$ws.Range("E42").Value = '  +7.07%  '
$ws.Range("E43").Value = '  +3.01%  '
$ws.Range("D44").Value = "'134.46"
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").Value = "'17.25"
$ws.Range("E45").Value = '  +2.36%  '
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").Value = "'3.91"
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").Value = "'22.27"
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").Value = '2.209.49'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").Value = '3.787.82'
$ws.Range("E51").Value = '  +5.68%  '
